$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.588.13'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.922.24'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4834'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.97%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2902'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06808'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '112.24'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.39%  '
$ws.Range("E11").Value = '  +5.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.916.96'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.488'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07577'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6738'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '294.54'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.565.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007688'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.520'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.160.03'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9990'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.460'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.502'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.32'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.105'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1068'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.53%  '
$ws.Range("E30").Value = '  +2.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.140'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.072'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04974'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7365'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.140'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02033'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.714'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.690'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.028'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '109.65'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4448'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8720'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.877'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.34%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.30%  '
$ws.Range("E46").Value = '  -0.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '49.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.261'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.78%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1232'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.07%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.90'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.42%  '
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2508'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.40%  '
